# Apply the "Edit Blog and Incomplete Multiple Browser Executions" update:
#  - Refresh the TEST-A-* / TEST-AS-* identifiers used across the three
#    product-management test sheets to the new sequence.
#  - Move the active/selected tab from ProductSheet (1st sheet) to
#    synchronization_sheet (3rd sheet).

$wb = $excel.ActiveWorkbook

$productSheet = $wb.Worksheets.Item("ProductSheet")
$editSheet = $wb.Worksheets.Item("editSheet")
$syncSheet = $wb.Worksheets.Item("synchronization_sheet")

# New "name" values for ProductSheet (column B) / editSheet (column A)
$newNames = @("TEST-A-68", "TEST-A-69", "TEST-A-70", "TEST-A-71", "TEST-A-72", "TEST-A-73")

for ($i = 0; $i -lt $newNames.Length; $i++) {
    $row = $i + 2
    $productSheet.Range("B$row").Value = $newNames[$i]
    $editSheet.Range("A$row").Value = $newNames[$i]
}

# New "name" values for synchronization_sheet (column B)
$newSyncNames = @("TEST-SA-56", "TEST-SA-57", "TEST-SA-58", "TEST-SA-59", "TEST-SA-60", "TEST-SA-61")

for ($i = 0; $i -lt $newSyncNames.Length; $i++) {
    $row = $i + 2
    $syncSheet.Range("B$row").Value = $newSyncNames[$i]
}

# Switch the active/selected tab from ProductSheet to synchronization_sheet
$syncSheet.Activate()
